# Updates crypto price/volume data to match the latest scrape.
# D-column "Price" values are plain text (e.g. "20.50", "1.01") that
# Excel would otherwise auto-convert to numbers (dropping trailing
# zeros); force text entry via NumberFormat "@", then restore the
# default "Normal" style so no stray per-cell style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.901.20'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.449.13'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.90%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '480.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +8.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +12.70%  '
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.504'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.462.71'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0968'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.48'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.327'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.08%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.122'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.894.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '55.039.00'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000135'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +13.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.481.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +10.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '315.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.994'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.167'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.406'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.613.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0776'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +15.28%  '
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.13'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.37%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.11'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +11.48%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.855'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.64'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.02'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.601'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0547'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.57%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.05%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.963.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0903'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '252.22'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +29.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0221'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.11%  '
